$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'92.482.87"
$ws.Range("E2").Value = "  +0.62%  "
$ws.Range("D3").Value = "'3.106.91"
$ws.Range("E3").Value = "  -0.69%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'234.78"
$ws.Range("E5").Value = "  -3.41%  "
$ws.Range("D6").Value = "'612.36"
$ws.Range("E6").Value = "  -1.04%  "
$ws.Range("E7").Value = "  -1.52%  "
$ws.Range("E8").Value = "  -0.28%  "
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("B10").Value = "LidoStakedEther"
$ws.Range("C10").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D10").Value = "'3.104.52"
$ws.Range("E10").Value = "  -0.71%  "
$ws.Range("B11").Value = "Cardano"
$ws.Range("C11").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D11").Value = "'0.790"
$ws.Range("E11").Value = "  +4.07%  "
$ws.Range("E12").Value = "  -3.88%  "
$ws.Range("D13").Value = "'0.0000243"
$ws.Range("E13").Value = "  -3.15%  "
$ws.Range("D14").Value = "'92.291.76"
$ws.Range("E14").Value = "  +0.76%  "
$ws.Range("D15").Value = "'33.84"
$ws.Range("E15").Value = "  -4.60%  "
$ws.Range("D16").Value = "'5.40"
$ws.Range("E16").Value = "  -3.84%  "
$ws.Range("D17").Value = "'3.688.95"
$ws.Range("E17").Value = "  -0.59%  "
$ws.Range("D18").Value = "'3.093.19"
$ws.Range("D19").Value = "'3.79"
$ws.Range("E19").Value = "  -0.25%  "
$ws.Range("D20").Value = "'14.50"
$ws.Range("E20").Value = "  -3.47%  "
$ws.Range("D21").Value = "'5.84"
$ws.Range("E21").Value = "  -1.60%  "
$ws.Range("E22").Value = "  +0.26%  "
$ws.Range("D23").Value = "'438.28"
$ws.Range("E23").Value = "  -4.25%  "
$ws.Range("D24").Value = "'9.11"
$ws.Range("E24").Value = "  -1.79%  "
$ws.Range("D25").Value = "'8.19"
$ws.Range("E25").Value = "  +4.98%  "
$ws.Range("D26").Value = "'5.57"
$ws.Range("E26").Value = "  -6.80%  "
$ws.Range("D27").Value = "'85.46"
$ws.Range("E27").Value = "  -4.36%  "
$ws.Range("D28").Value = "'11.47"
$ws.Range("E28").Value = "  -2.28%  "
$ws.Range("D29").Value = "'3.272.81"
$ws.Range("D31").Value = "'0.181"
$ws.Range("E31").Value = "  +7.42%  "
$ws.Range("E32").Value = "  +3.60%  "
$ws.Range("E33").Value = "  -15.41%  "
$ws.Range("E34").Value = "  -29.44%  "
$ws.Range("D35").Value = "'9.16"
$ws.Range("E35").Value = "  -3.12%  "
$ws.Range("D36").Value = "'8.08"
$ws.Range("E36").Value = "  +7.54%  "
$ws.Range("D37").Value = "'0.163"
$ws.Range("E37").Value = "  -7.85%  "
$ws.Range("D38").Value = "'25.62"
$ws.Range("E38").Value = "  -3.43%  "
$ws.Range("D39").Value = "'3.97"
$ws.Range("E39").Value = "  +2.63%  "
$ws.Range("E40").Value = "  -3.49%  "
$ws.Range("E41").Value = "  +7.61%  "
$ws.Range("D43").Value = "'463.27"
$ws.Range("E43").Value = "  -5.78%  "
$ws.Range("D44").Value = "'0.427"
$ws.Range("E44").Value = "  -3.32%  "
$ws.Range("D45").Value = "'3.30"
$ws.Range("E45").Value = "  -2.72%  "
$ws.Range("E46").Value = "  +0.03%  "
$ws.Range("D47").Value = "'159.89"
$ws.Range("E47").Value = "  +2.14%  "
$ws.Range("D48").Value = "'0.680"
$ws.Range("E48").Value = "  -4.32%  "
$ws.Range("D50").Value = "'1.32"
$ws.Range("E50").Value = "  -2.38%  "
$ws.Range("B51").Value = "OKB"
$ws.Range("C51").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D51").Value = "'43.77"
$ws.Range("E51").Value = "  -0.60%  "

# The leading apostrophe marks these as "text quote prefix" cells in Excel's
# model; clear that formatting flag so the cells end up as plain text cells
# (same as the rest of the sheet) rather than carrying a quote-prefix style.
$ws.Range("D2:D51").ClearFormats()
